# Insert a new data row before the current row 27 (shifting the existing
# rows 27-34 down to 28-35) and populate the new row with the latest
# weekly price observation for "Feria Lagunitas de Puerto Montt - Alcachofa".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 27:34 down to 28:35, opening up a blank row 27.
$ws.Rows("27:27").Insert()

# Fill the newly-opened row 27 with the new weekly record.
$ws.Range("A27").Value2 = 4
$ws.Range("B27").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C27").Value2 = "Los Lagos"
$ws.Range("D27").Value2 = 44855
$ws.Range("E27").Value2 = 10
$ws.Range("F27").Value2 = 100112013
$ws.Range("G27").Value2 = "Alcachofa"
$ws.Range("H27").Value2 = "Madrigal"
$ws.Range("I27").Value2 = "Primera"
$ws.Range("J27").Value2 = 160
$ws.Range("K27").Value2 = 10000
$ws.Range("L27").Value2 = 10000
$ws.Range("M27").Value2 = 10000
$ws.Range("N27").Value2 = "`$/caja 30 unidades"
$ws.Range("O27").Value2 = "Provincia de Limarí"
$ws.Range("P27").Value2 = 333
$ws.Range("Q27").Value2 = 30
$ws.Range("R27").Value2 = "Hortaliza"
